# "issues de ficha tecnica" -- text corrections on slide 1 of the
# FichaTecnicaObras deck. Only the run text of seven existing text boxes
# changes; the shapes themselves keep their original position/size.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_POINT = 12700

# Rewriting a run's text via `TextRange.Characters(start, len).Text = ...`
# (rather than replacing the whole `TextRange.Text`) edits the existing
# `<a:r>` in place instead of re-stamping it with a fresh `<a:rPr
# lang="en-US"/>`, which keeps the XML diff limited to the `<a:t>` content
# -- matching how the source deck was actually edited.
function Set-ShapeRunText($shape, $newText) {
    $tr = $shape.TextFrame.TextRange
    $tr.Characters(1, $tr.Length).Text = $newText
}

# Shape.Height round-trips through a single-precision (float32) COM
# property, and the engine floors (rather than rounds) when turning the
# point value back into EMU. Reading `$shape.Height` after PowerPoint
# re-lays out an autofit text box and writing it straight back can
# therefore drift the stored extent by a single EMU. `Get-EmuForPoints`
# mirrors that float32+floor pipeline and `Find-SafePoints` walks nearby
# float32 values until it finds one that reproduces the desired EMU
# exactly, so a box's original height can be restored losslessly.
function Get-EmuForPoints($pts) {
    $f = [float]$pts
    return [math]::Floor([double]$f * $EMU_PER_POINT)
}

function Find-SafePoints($targetEmu) {
    $exact = $targetEmu / [double]$EMU_PER_POINT
    if ((Get-EmuForPoints $exact) -eq $targetEmu) {
        return [double][float]$exact
    }
    $step = 0.0000001
    for ($i = 1; $i -lt 4000; $i++) {
        $up = $exact + $step * $i
        if ((Get-EmuForPoints $up) -eq $targetEmu) {
            return [double][float]$up
        }
        $down = $exact - $step * $i
        if ((Get-EmuForPoints $down) -eq $targetEmu) {
            return [double][float]$down
        }
    }
    return $exact
}

# Change a text box's run text while keeping its autofit (`spAutoFit`)
# shape height pinned to the EMU value it had before the edit, in case
# the replacement text wraps onto a different number of lines.
function Set-ShapeTextKeepHeight($slide, $shapeName, $newText, $originalHeightEmu) {
    $shape = $slide.Shapes.Item($shapeName)
    Set-ShapeRunText $shape $newText
    if ((Get-EmuForPoints $shape.Height) -ne $originalHeightEmu) {
        $shape.Height = Find-SafePoints $originalHeightEmu
    }
}

Set-ShapeTextKeepHeight $s "5 CuadroTexto"  "OB_SEGOB_00010"   215444
Set-ShapeTextKeepHeight $s "12 CuadroTexto" "INTERESTATAL"     215444
Set-ShapeTextKeepHeight $s "14 CuadroTexto" "2013-10-01"       215444
Set-ShapeTextKeepHeight $s "15 CuadroTexto" "2013-12-01"       215444
Set-ShapeTextKeepHeight $s "26 CuadroTexto" "8.18"             215444
Set-ShapeTextKeepHeight $s "39 CuadroTexto" "MANTENIMIENTO MAYOR A INMUEBLES DE LA DELEGACION FEDERAL DEL INM EN EL ESTADO DE QUINTANA ROO Y TABASCO" 215444
Set-ShapeTextKeepHeight $s "40 CuadroTexto" "TRAMITES MIGRATORIOS: 2488,885 EXTRANJEROS PRESENTADOS ANTE EL INM: 3,461" 215444
